$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing data rows (2..185) down by one, from the bottom up,
# so we don't clobber values while moving them. This avoids a whole-row
# Insert(), which would drag header formatting/extra blank cells onto the
# freshly inserted row.
for ($r = 185; $r -ge 2; $r--) {
    $ws.Cells.Item($r + 1, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 3).Value2 = $ws.Cells.Item($r, 3).Value2
}

# Populate the newly freed row 2 with the new song entry
$ws.Cells.Item(2, 2).Value2 = "Lana Del Rey - Henry, come on (Lyrics)"
$ws.Cells.Item(2, 3).Value2 = "https://www.youtube.com/watch?v=wasFuIuPh5k"
